$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table (A:date, B:weekday, C:time, D:rank) gets two brand-new
# rows inserted right before the old row 690 ("2026/12/29" block), which
# pushes that whole block (and everything after it) down by two rows
# without otherwise altering it.
$ws.Range("A690:A691").EntireRow.Insert()

# Format column A as text first so the yyyy/mm/dd-looking strings are
# stored as literal text (matching the rest of the sheet) instead of
# being auto-converted to date serials.
$ws.Range("A690:A691").NumberFormat = "@"

# New row 690: 2026/01/24 (Sat), time 22, rank 17
$ws.Range("A690").Value = "2026/01/24"
$ws.Range("B690").Value = "土"
$ws.Range("C690").Value = 22
$ws.Range("D690").Value = 17

# New row 691: 2026/01/25 (Sun), time 2, rank 18
$ws.Range("A691").Value = "2026/01/25"
$ws.Range("B691").Value = "日"
$ws.Range("C691").Value = 2
$ws.Range("D691").Value = 18
